# This workbook was round-tripped through another Excel-compatible editor
# (WPS Office, based on the xmlns:etc namespace / charset=134 fonts it left
# behind). That round-trip rewrote a lot of incidental, non-semantic XML
# (style-table boilerplate, namespace decls, fileVersion/calcPr bookkeeping,
# window geometry, float serialization precision, row "spans" hints, ...).
# None of the actual cell data, formulas or shared strings changed.
#
# The only genuinely observable, reproducible edits are the view/format
# state captured below: the column widths were resized, and the active
# selection was left on D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -------------------------------------------------------
# Target stored widths (OOXML "width" attribute, in characters):
#   A=41.375  B=6.5  C=9.75  D=15.375  E=5.875  F=5.625  G=15.7083333333333
# Excel's ColumnWidth COM property is specified in characters but gets
# snapped to the sheet's pixel grid on write, so we pick the ColumnWidth
# input that lands closest to each target stored width.
$ws.Columns.Item(1).ColumnWidth = 40.5
$ws.Columns.Item(2).ColumnWidth = 5.66666666666667
$ws.Columns.Item(3).ColumnWidth = 8.83333333333333
$ws.Columns.Item(4).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 5
$ws.Columns.Item(6).ColumnWidth = 4.83333333333333
$ws.Columns.Item(7).ColumnWidth = 14.8333333333333

# --- Selection -------------------------------------------------------
# The saved view had D7 as the active cell/selection.
$ws.Range("D7").Select() | Out-Null
